$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('A1').Value = 'basketball with design'
$ws.Range('A2').Value = 'yoga knee'
$ws.Range('A3').Value = 'basketball skins'
$ws.Range('A4').Value = 'leg pads hockey'
$ws.Range('A5').Value = 'elastic calf sleeve'
$ws.Range('A6').Value = 'softball shorts girls'
$ws.Range('A7').Value = 'compression knee sleeve with pads'
$ws.Range('A8').Value = 'volleyball compression shorts'
$ws.Range('A9').Value = 'knee compression sleeve for basketball'
$ws.Range('A10').Value = 'youth compression sleeve baseball'
$ws.Range('A11').Value = 'youth basketball shorts'
$ws.Range('A12').Value = 'hockey pants'
$ws.Range('A13').Value = 'mens long compression shorts'
$ws.Range('A14').Value = 'basketball knee sleeve'
$ws.Range('A15').Value = 'running pants for men'
$ws.Range('A16').Value = 'tall mens tights'
$ws.Range('A17').Value = 'baseball compression'
$ws.Range('A18').Value = 'protective pad'
$ws.Range('A19').Value = 'baseball youth pants'
$ws.Range('A20').Value = 'kneeling pad for exercise'
$ws.Range('A21').Value = 'basketball compression sleeve youth'
$ws.Range('A22').Value = 'volleyball pads for girls'
$ws.Range('A23').Value = 'sheer protectors'
$ws.Range('A24').Value = 'compression pants with pouch'
$ws.Range('A25').Value = 'protect knee'
$ws.Range('A26').Value = 'calf sleeve weight'
$ws.Range('A27').Value = 'mens cycling pants with padding'
$ws.Range('A28').Value = 'hip pads sports'
$ws.Range('A29').Value = 'baseball shorts'
$ws.Range('A30').Value = 'air knee pads'
$ws.Range('A31').Value = 'knee pads work construction'
$ws.Range('A32').Value = 'knee sleeves for basketball pair'
$ws.Range('A33').Value = 'compression tights men pack'
$ws.Range('A34').Value = 'football leggings for boys'
$ws.Range('A35').Value = 'construction pants for men'
$ws.Range('A36').Value = 'knee compression sleeve youth'
$ws.Range('A37').Value = 'knee pads thin'
$ws.Range('A38').Value = 'mens leggings pouch'
$ws.Range('A39').Value = 'running pants men tights'
$ws.Range('A40').Value = 'patella band'
$ws.Range('A41').Value = 'youth knee sleeve'
$ws.Range('A42').Value = 'calf compression sleeve spandex'
$ws.Range('A43').Value = 'youth baseball compression sleeve'
$ws.Range('A44').Value = 'protective pads'
$ws.Range('A45').Value = 'mens leggings cold'
$ws.Range('A46').Value = 'youth boys compression'
$ws.Range('A47').Value = 'little boys baseball pants size 5'
$ws.Range('A48').Value = 'mens running tights pants'
$ws.Range('A49').Value = 'men tights sports'
$ws.Range('A50').Value = 'snowboarding padded shorts'
$ws.Range('A51').Value = 'baseball dirt'
$ws.Range('A52').Value = 'boy compression pants'
$ws.Range('A53').Value = 'youth kneepads'
$ws.Range('A54').Value = 'best construction knee pads'
$ws.Range('A55').Value = 'fall cycling pants'
$ws.Range('A56').Value = 'knee pads mountain biking'
$ws.Range('A57').Value = 'football hip pads'
$ws.Range('A58').Value = 'catchers leg guards adult'
$ws.Range('A59').Value = 'flexible knee pad'
$ws.Range('A60').Value = 'basketball shorts pack of 5'
$ws.Range('A61').Value = 'volleyball mens shorts'
$ws.Range('A62').Value = 'knee pads girls'
$ws.Range('A63').Value = 'bees knees'
$ws.Range('A64').Value = 'sport leggings'
$ws.Range('A65').Value = 'girls lacrosse shorts'
$ws.Range('A66').Value = 'best knee pads'
$ws.Range('A67').Value = 'calf pads'
$ws.Range('A68').Value = 'burns baseball'
$ws.Range('A69').Value = 'boys workout leggings'
$ws.Range('A70').Value = '6 short pants'
$ws.Range('A71').Value = 'basketball leg sleeve youth'
$ws.Range('A72').Value = 'calf compression leggings'
$ws.Range('A73').Value = 'capri leggings with mesh'
$ws.Range('A74').Value = 'basketball sleeve for youth'
$ws.Range('A75').Value = 'extra thick knee pads'
$ws.Range('A76').Value = 'girdle football adult'
$ws.Range('A77').Value = 'indoor pants'
$ws.Range('A78').Value = 'youth girls knee pads'
$ws.Range('A79').Value = 'mens volleyball knee sleeve'
$ws.Range('A80').Value = 'knees for men'
$ws.Range('A81').Value = 'patella band knee'
$ws.Range('A82').Value = 'football knee sleeves'
$ws.Range('A83').Value = 'basketball sleeve for men'
$ws.Range('A84').Value = 'paintball pads and protection'
$ws.Range('A85').Value = 'sports leggings boys'
$ws.Range('A86').Value = 'football calf sleeve'
$ws.Range('A87').Value = '$5 and below'
$ws.Range('A88').Value = 'capri tights'
$ws.Range('A89').Value = 'compression bands for knees'
$ws.Range('A90').Value = 'mens compression pants black'
$ws.Range('A91').Value = 'playing ball on running water'
$ws.Range('A92').Value = 'youth large baseball pants'
$ws.Range('A93').Value = 'girls spandex shorts black volleyball'
$ws.Range('A94').Value = 'football short tights'
$ws.Range('A95').Value = 'thigh pads'
$ws.Range('A96').Value = 'tights and leggings'
$ws.Range('A97').Value = 'professional knee pad'
$ws.Range('A98').Value = 'short youth baseball pants'
$ws.Range('A99').Value = 'youth calf compression sleeve'
$ws.Range('A100').Value = 'knee guards'
